# "harmonized similar tags to be the same"
#
# On the "isa_template" sheet, the ER (entity role) row (row 12) and its
# accompanying "ER Term Accession Number" row (row 13) used free-text /
# inconsistent casing ("Mass spectrometry", "MS", "Data", "Processing")
# instead of reusing the harmonized ontology term + accession already used
# elsewhere in the sheet ("Mass Spectrometry" / NCIT:C17156 and
# "data processing" / NCIT:C47925).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("isa_template")

# Row 12: #ER list -> ER: "Data" becomes "data processing", "Mass spectrometry"
# becomes "Mass Spectrometry" (harmonized casing), the now-redundant "MS" and
# "Processing" tag cells are cleared.
$ws.Range("C12").Value = "data processing"
$ws.Range("D12").Value = "Mass Spectrometry"
$ws.Range("E12").ClearContents()
$ws.Range("F12").ClearContents()

# Row 13: ER Term Accession Number - fill in the matching NCIT accession
# numbers for the two harmonized terms above, and grow the row to fit the
# wrapped text.
$ws.Range("C13").Value = "NCIT:C47925"
$ws.Range("D13").Value = "NCIT:C17156"
$ws.Rows.Item(13).RowHeight = 43.2

# Reflect the editor's last selection on the sheet.
$ws.Range("G12").Select()
